# [ADD] new parameters set
# Update the RS11 parameter table with a new set of values and extend the
# sheet with a few blank (underlined-style) rows below the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing parameter values -------------------------------------
# A2 "m"   (vehicle mass)
$ws.Range("B2").Value = 350
# A3 "l_p" (pedal - cylinder ratio)
$ws.Range("B3").Value = 4.33
# A4 "Amc" (master cylinder cross section area) - same value in both columns
$ws.Range("B4").Value = 1.98
$ws.Range("C4").Value = 1.98
# A8 "R"   (external wheel radius) - same value in both columns
$ws.Range("B8").Value = 230
$ws.Range("C8").Value = 230

# --- Extend the sheet with three new (empty) underlined cells -------------
# Mirrors the look of the explanatory-text column (E) being carried a few
# rows further down, growing the used range to A1:E11.
$ws.Range("E9").Font.Underline = $true
$ws.Range("E10").Font.Underline = $true
$ws.Range("E11").Font.Underline = $true

# --- Page setup -------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection --------------------------------------------------------------
$ws.Range("B2").Select() | Out-Null
